# Update the Fitness column (C) for rows 2 through 252 to the new value 7310.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C252").Value = 7310
